$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q4" right before the existing
#    "2022-Q3" sheet (which currently sits right after the "总计" sheet).
# ---------------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($anchor)
$newSheet.Name = "2022-Q4"

# Clone the look & feel (header row style + first-column style) from the
# neighbouring quarter sheet so the new sheet matches the existing ones.
$template = $wb.Worksheets.Item("2022-Q3")
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A16").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows (fund holdings for 2022-Q4)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'008850"
$newSheet.Range("C2").Value = "景顺长城价值稳进三年定期开放灵活配置混合"
$newSheet.Range("D2").Value = "'18.18"
$newSheet.Range("E2").Value = "'97.45"
$newSheet.Range("F2").Value = "'3.94"
$newSheet.Range("G2").Value = "'0.7163"
$newSheet.Range("H2").Value = 9
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'009846"
$newSheet.Range("C3").Value = "富兰克林国海港股通远见价值混合"
$newSheet.Range("D3").Value = "'15.17"
$newSheet.Range("E3").Value = "'90.16"
$newSheet.Range("F3").Value = "'3.17"
$newSheet.Range("G3").Value = "'0.4809"
$newSheet.Range("H3").Value = 7
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'009098"
$newSheet.Range("C4").Value = "景顺长城价值领航两年持有期混合"
$newSheet.Range("D4").Value = "'7.08"
$newSheet.Range("E4").Value = "'92.98"
$newSheet.Range("F4").Value = "'4.02"
$newSheet.Range("G4").Value = "'0.2846"
$newSheet.Range("H4").Value = 8
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'008715"
$newSheet.Range("C5").Value = "景顺长城价值驱动一年持有期灵活配置混合"
$newSheet.Range("D5").Value = "'6.75"
$newSheet.Range("E5").Value = "'92.52"
$newSheet.Range("F5").Value = "'4.15"
$newSheet.Range("G5").Value = "'0.2801"
$newSheet.Range("H5").Value = 6
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'008060"
$newSheet.Range("C6").Value = "景顺长城价值边际灵活配置混合A"
$newSheet.Range("D6").Value = "'6.42"
$newSheet.Range("E6").Value = "'91.73"
$newSheet.Range("F6").Value = "'3.70"
$newSheet.Range("G6").Value = "'0.2375"
$newSheet.Range("H6").Value = 8
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'015779"
$newSheet.Range("C7").Value = "景顺长城价值边际灵活配置混合C"
$newSheet.Range("D7").Value = "'2.40"
$newSheet.Range("E7").Value = "'91.73"
$newSheet.Range("F7").Value = "'3.70"
$newSheet.Range("G7").Value = "'0.0888"
$newSheet.Range("H7").Value = 8
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'004497"
$newSheet.Range("C8").Value = "前海开源多元策略灵活配置混合C"
$newSheet.Range("D8").Value = "'1.79"
$newSheet.Range("E8").Value = "'79.66"
$newSheet.Range("F8").Value = "'4.54"
$newSheet.Range("G8").Value = "'0.0813"
$newSheet.Range("H8").Value = 3
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "'004496"
$newSheet.Range("C9").Value = "前海开源多元策略灵活配置混合A"
$newSheet.Range("D9").Value = "'1.30"
$newSheet.Range("E9").Value = "'79.66"
$newSheet.Range("F9").Value = "'4.54"
$newSheet.Range("G9").Value = "'0.0590"
$newSheet.Range("H9").Value = 3
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "'011471"
$newSheet.Range("C10").Value = "鹏华致远成长混合A"
$newSheet.Range("D10").Value = "'1.67"
$newSheet.Range("E10").Value = "'60.84"
$newSheet.Range("F10").Value = "'2.22"
$newSheet.Range("G10").Value = "'0.0371"
$newSheet.Range("H10").Value = 1
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "'004098"
$newSheet.Range("C11").Value = "前海开源港股通股息率50强股票"
$newSheet.Range("D11").Value = "'0.53"
$newSheet.Range("E11").Value = "'90.79"
$newSheet.Range("F11").Value = "'6.47"
$newSheet.Range("G11").Value = "'0.0343"
$newSheet.Range("H11").Value = 1
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "'006923"
$newSheet.Range("C12").Value = "前海开源沪港深非周期性行业股票A"
$newSheet.Range("D12").Value = "'0.28"
$newSheet.Range("E12").Value = "'90.65"
$newSheet.Range("F12").Value = "'4.54"
$newSheet.Range("G12").Value = "'0.0127"
$newSheet.Range("H12").Value = 8
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "'006924"
$newSheet.Range("C13").Value = "前海开源沪港深非周期性行业股票C"
$newSheet.Range("D13").Value = "'0.24"
$newSheet.Range("E13").Value = "'90.65"
$newSheet.Range("F13").Value = "'4.54"
$newSheet.Range("G13").Value = "'0.0109"
$newSheet.Range("H13").Value = 8
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "'161124"
$newSheet.Range("C14").Value = "易方达香港恒生综合小型股指数（QDII-LOF）A"
$newSheet.Range("D14").Value = "'0.24"
$newSheet.Range("E14").Value = "'94.45"
$newSheet.Range("F14").Value = "'1.53"
$newSheet.Range("G14").Value = "'0.0037"
$newSheet.Range("H14").Value = 4
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "'011472"
$newSheet.Range("C15").Value = "鹏华致远成长混合C"
$newSheet.Range("D15").Value = "'0.06"
$newSheet.Range("E15").Value = "'60.84"
$newSheet.Range("F15").Value = "'2.22"
$newSheet.Range("G15").Value = "'0.0013"
$newSheet.Range("H15").Value = 1
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "'006263"
$newSheet.Range("C16").Value = "易方达香港恒生综合小型股指数（QDII-LOF）C"
$newSheet.Range("D16").Value = "'0.05"
$newSheet.Range("E16").Value = "'94.45"
$newSheet.Range("F16").Value = "'1.53"
$newSheet.Range("G16").Value = "'0.0008"
$newSheet.Range("H16").Value = 4

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new top row for 2022-Q4 and
#    shift the existing quarterly summary rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Give the brand-new row 8 the same look as the existing numbered rows
# before writing into it (copy formatting from row 7's A cell).
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)

# Fill bottom-up so each write lands on a cell that already owns the
# correct pre-existing style.
$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 1
$summary.Range("D8").Value = 0.18

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 1
$summary.Range("D7").Value = 0.14

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q4"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 0.05

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.09

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 4
$summary.Range("D4").Value = 0.25

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 8
$summary.Range("D3").Value = 0.17

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 15
$summary.Range("D2").Value = 2.33
